$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - August (through ...) updates
$ws.Range("A9").Value = "August (through 08-31)"
$ws.Range("B9").Value = 32
$ws.Range("C9").Value = 79
$ws.Range("F9").Value = 45
$ws.Range("G9").Value = 163
$ws.Range("H9").Value = 156

# Row 10 - Total updates
$ws.Range("B10").Value = 194
$ws.Range("C10").Value = 381
$ws.Range("F10").Value = 349
$ws.Range("G10").Value = 784
$ws.Range("H10").Value = 1069
